$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rearrange the columns so userId is first, then title, then body
$data = @(
    @("userId", "title", "body"),
    @(1, "First Post -DDT", "This is the body of the first post created using DDT-excel with Playwright API tests."),
    @(2, "Second Post- DDT", "This is the body of the second post created using DDT-excel with Playwright API tests."),
    @(3, "Thrid Post- DDT", "This is the thrid post created via DDT-excel.")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    for ($c = 0; $c -lt 3; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $data[$r][$c]
    }
}

[void]$ws.Range("A1:C4").Select()
